# github.docx interop edit: add a French (fr-FR) set of GitHub/GitHub-Pages
# API token lines plus the Guadeloupe bus-stop/coordinate list below the token,
# reproducing the paragraph-mark language tagging Word's French proofing pass
# applies as each line is typed (including the spell-check proofErr wrapping
# around the unrecognised proper nouns and the lastRenderedPageBreak Word left
# mid-session).

$d = $word.ActiveDocument

# 1) The pre-existing token paragraph only picked up the French paragraph-mark
#    language (no text/run changed), so splice just a <w:pPr> in via InsertXML
#    while leaving the original run untouched.
$p1 = $d.Paragraphs(1)
$p1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="67809813" w14:textId="5EDA266B" w:rsidR="00C271C8" w:rsidRDefault="00C271C8"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r w:rsidRPr="00C271C8"><w:t>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</w:t></w:r></w:p>
'@
$p1.Range.InsertXML($p1Xml)

# 2) Append each new fr-FR paragraph (second token, then the alternating
#    place-name / "lat, long" lines) at the end of the document, in order.
#    Re-fetch the end-of-document range each time so successive inserts land
#    after the previously inserted paragraph instead of all at one stale spot.

# paragraph 2
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>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</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 3
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Saint-François</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 4
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.250453094513063, -61.272616319469535</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 5
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Darboussier</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 6
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.23482950072778, -61.53309884478915</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 7
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Périnet</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 8
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.213118094624487, -61.47718754941982</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 9
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Grande Ravine</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 10
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.215967972649068, -61.47961760291421</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 11
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Riviéra</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 12
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.217807997101172, -61.485820926274876</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 13
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Carr La Bouaye</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 14
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.217824741851715, -61.491074110863245</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 15
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Belle Plaine</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 16
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.217613179304927, -61.49451872882002</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 17
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Ecole de Poucet</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 18
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.219060659497853, -61.502375691746224</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 19
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Grand Baie</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 20
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.218310277388863, -61.511290241477454</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 21
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Bas du Fort</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 22
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.21976324630858, -61.51717010643815</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 23
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Blanchard</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 24
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.223477702406306, -61.52514039639234</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 25
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Blanchard Université</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 26
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.224759430691957, -61.52707538890146</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 27
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Université </w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 28
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:lastRenderedPageBreak/><w:t>16.2269500838484, -61.52927029877085</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 29
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Louis </w:t></w:r><w:proofErr w:type="spellStart"/><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Douldat</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 30
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.230504637063415, -61.5313383785461</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 31
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Raspail</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 32
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.234273746919115, -61.53177015373589</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 33
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>Camille Desmoulin</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

# paragraph 34
$end = $d.Content
$end.Collapse(0)
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r>
        <w:rPr>
          <w:lang w:val="fr-FR"/></w:rPr><w:t>16.23580345914415, -61.533023844685694</w:t></w:r></w:p>
'@
$end.InsertXML($paraXml)

